# Lunchtimes and lecturer time limits implemented.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Replace the lecturer in column B: "Melina Aldinger" -> "Jane Courtney" ---
$ws1.Range("B2").Value = "Jane Courtney "

# --- Fill the lecturer-expertise grid (B3:U14) with the new checkerboard pattern ---
for ($r = 3; $r -le 14; $r++) {
    for ($c = 2; $c -le 21; $c++) {
        if ((($r + $c) % 2) -eq 1) {
            $val = 1
        } else {
            $val = 0
        }
        $ws1.Cells.Item($r, $c).Value = $val
    }
}

# --- Update the visible selection / scroll position on Sheet1 back to A1:U14 ---
$ws1.Range("A1:U14").Select()

# --- Add a new, empty "Sheet2" right after "Sheet1" ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"
$ws2.Range("A1:U14").Select()

# --- Re-activate Sheet1 so it remains the selected/visible tab ---
$ws1.Activate()
